$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v0 = @'
{
  "script": "vulnerable",
  "score": 5,
  "findings": [
    {
      "line": 2,
      "severity": "Warning",
      "statement": "$filespath = '##win.monitored.files##';",
      "reason": "Sensitive variable declared but not used in a risky sink.",
      "recommendation": "Use this variable only for internal authentication purposes within the script; do not expose its value externally.",
      "code_suggestion": "$filespath = 'internal-use-only'"
    },
    {
      "line": 3,
      "severity": "Warning",
      "statement": "$hostname = '##system.hostname##';",
      "reason": "Sensitive variable declared but not used in a risky sink.",
      "recommendation": "Use this variable only for internal authentication purposes within the script; do not expose its value externally.",
      "code_suggestion": "$hostname = 'internal-use-only'"
    },
    {
      "line": 4,
      "severity": "Warning",
      "statement": "$user = '##wmi.user##';",
      "reason": "Sensitive variable declared but not used in a risky sink.",
      "recommendation": "Use this variable only for internal authentication purposes within the script; do not expose its value externally.",
      "code_suggestion": "$user = 'internal-use-only'"
    },
    {
      "line": 5,
      "severity": "Error",
      "statement": "$pass = '##wmi.pass##';",
      "reason": "Sensitive variable is used in a risky sink (file output).",
      "recommendation": "Store this credential securely and only use it for authentication within the script; never write it to files.",
      "code_suggestion": "$pass = New-Object PSCredential('', 'internal-use-only')"
    },
    {
      "line": 8,
      "severity": "Error",
      "statement": "Set-Content -Path \"sensitive.txt\" -Value $nonSensitive",
      "reason": "Sensitive variable is written to a file (leaked).",
      "recommendation": "Never write sensitive credentials to files; use them only for internal authentication.",
      "code_suggestion": "Remove this line or replace with safe operation"
    },
    {
      "line": 14,
      "severity": "Error",
      "statement": "Set-Content -Path \"sensitive.txt\" -Value $testVar",
      "reason": "Sensitive variable is written to a file (leaked).",
      "recommendation": "Never write sensitive credentials to files; use them only for internal authentication.",
      "code_suggestion": "Remove this line or replace with safe operation"
    },
    {
      "line": 20,
      "severity": "Error",
      "statement": "$newVar = $testVar1",
      "reason": "Sensitive variable is copied and then potentially leaked.",
      "recommendation": "Avoid copying sensitive variables; use them directly for authentication only.",
      "code_suggestion": "Remove this line or replace with direct usage"
    },
    {
      "line": 30,
      "severity": "Error",
      "statement": "Set-Content -Path \"sensitive.txt\" -Value $new1Var",
      "reason": "Sensitive variable is written to a file (leaked).",
      "recommendation": "Never write sensitive credentials to files; use them only for internal authentication.",
      "code_suggestion": "Remove this line or replace with safe operation"
    }
  ]
}
'@
$ws.Range("B2").Value = $v0

$v1 = @'
{
  "script": "safe",
  "score": 10,
  "findings": [
    {
      "line": 3,
      "severity": "Warning",
      "statement": "$password = \"YourSecurePassword123\"",
      "reason": "Sensitive variable declared but not used in a risky sink.",
      "recommendation": "This password should be used strictly for authentication within the script and never printed or logged.",
      "code_suggestion": "$password = Read-Host -AsSecureString -Prompt \"Enter password:\""
    },
    {
      "line": 4,
      "severity": "Warning",
      "statement": "$securePassword = ConvertTo-SecureString $password -AsPlainText -Force",
      "reason": "Sensitive variable converted to secure string but not used in a risky sink.",
      "recommendation": "The secure password should be used only for authentication and never printed or logged.",
      "code_suggestion": "$securePassword = Read-Host -AsSecureString -Prompt \"Enter password:\""
    },
    {
      "line": 5,
      "severity": "Warning",
      "statement": "$cred = New-Object System.Management.Automation.PSCredential($username, $securePassword)",
      "reason": "Sensitive credentials object created but not used in a risky sink.",
      "recommendation": "Use this credential only for authentication and never print or log it.",
      "code_suggestion": "$cred = Get-Credential -UserName $username -Message \"Enter password:\""
    },
    {
      "line": 17,
      "severity": "Warning",
      "statement": "Write-Host \"Hostname: $($cs.Name)\"",
      "reason": "Sensitive hostname information is being printed to the console.",
      "recommendation": "Only print this if absolutely necessary for debugging or auditing purposes; otherwise, remove the output.",
      "code_suggestion": "# Write-Host \"Hostname: $($cs.Name)\""
    },
    {
      "line": 18,
      "severity": "Warning",
      "statement": "Write-Host \"OS: $($os.Caption) ($($os.OSArchitecture))\"",
      "reason": "Sensitive OS information is being printed to the console.",
      "recommendation": "Only print this if absolutely necessary; otherwise, remove the output.",
      "code_suggestion": "# Write-Host \"OS: $($os.Caption) ($($os.OSArchitecture))\""
    },
    {
      "line": 19,
      "severity": "Warning",
      "statement": "Write-Host \"CPU: $($cpu.Name)\"",
      "reason": "Sensitive CPU information is being printed to the console.",
      "recommendation": "Only print this if absolutely necessary; otherwise, remove the output.",
      "code_suggestion": "# Write-Host \"CPU: $($cpu.Name)\""
    },
    {
      "line": 20,
      "severity": "Warning",
      "statement": "Write-Host \"Total Physical Memory: $([math]::Round($cs.TotalPhysicalMemory / 1GB, 2)) GB\"",
      "reason": "Sensitive memory information is being printed to the console.",
      "recommendation": "Only print this if absolutely necessary; otherwise, remove the output.",
      "code_suggestion": "# Write-Host \"Total Physical Memory: $([math]::Round($cs.TotalPhysicalMemory / 1GB, 2)) GB\""
    },
    {
      "line": 24,
      "severity": "Warning",
      "statement": "$procs | Sort-Object -Property ProcessId | Select-Object -First 5",
      "reason": "Process information is being printed to the console.",
      "recommendation": "Only print this if absolutely necessary for debugging or auditing purposes; otherwise, remove the output.",
      "code_suggestion": "# $procs | Sort-Object -Property ProcessId | Select-Object -First 5"
    }
  ]
}
'@
$ws.Range("B3").Value = $v1

$v2 = @'
{
  "script": "vulnerable",
  "score": 3,
  "findings": [
    {
      "line": 3,
      "severity": "Warning",
      "statement": "$password = \"YourSecurePassword123!\"",
      "reason": "Sensitive password stored in plain text variable.",
      "recommendation": "Use a secure method to store or retrieve credentials. Avoid hardcoding passwords directly in scripts.",
      "code_suggestion": "$password = Read-Host -AsSecureString 'Enter password:'"
    },
    {
      "line": 13,
      "severity": "Error",
      "statement": "Write-Host \"[+] Gathering system info from $remoteComputer using secure WMI session...\" -ForegroundColor Cyan",
      "reason": "Sensitive variable $remoteComputer is being printed to the console.",
      "recommendation": "Remove or mask sensitive information before printing. Only display necessary details for authorized users.",
      "code_suggestion": "Write-Host \"[+] Gathering system info from a remote computer using secure WMI session...\" -ForegroundColor Cyan"
    },
    {
      "line": 19,
      "severity": "Error",
      "statement": "Write-Host \"Hostname         : $($cs.Name)\"",
      "reason": "Sensitive variable $cs (computer system object) is being printed to the console.",
      "recommendation": "Remove or mask sensitive information before printing. Only display necessary details for authorized users.",
      "code_suggestion": "Write-Host \"Hostname         : [Computer Name Masked]\""
    },
    {
      "line": 20,
      "severity": "Error",
      "statement": "Write-Host \"OS               : $($os.Caption) ($($os.OSArchitecture))\"",
      "reason": "Sensitive variable $os (operating system object) is being printed to the console.",
      "recommendation": "Remove or mask sensitive information before printing. Only display necessary details for authorized users.",
      "code_suggestion": "Write-Host \"OS               : [Operating System Masked]\""
    },
    {
      "line": 21,
      "severity": "Error",
      "statement": "Write-Host \"CPU              : $($cpu.Name)\"",
      "reason": "Sensitive variable $cpu (processor object) is being printed to the console.",
      "recommendation": "Remove or mask sensitive information before printing. Only display necessary details for authorized users.",
      "code_suggestion": "Write-Host \"CPU              : [Processor Name Masked]\""
    },
    {
      "line": 22,
      "severity": "Error",
      "statement": "Write-Host \"Total Memory     : $([math]::Round($cs.TotalPhysicalMemory / 1GB, 2)) GB\"",
      "reason": "Sensitive variable $cs (computer system object) is being printed to the console.",
      "recommendation": "Remove or mask sensitive information before printing. Only display necessary details for authorized users.",
      "code_suggestion": "Write-Host \"Total Memory     : [Memory Size Masked] GB\""
    },
    {
      "line": 25,
      "severity": "Error",
      "statement": "Write-Host \"`nDisks:\n\"",
      "reason": "Sensitive variable $disks (disk objects) is being printed to the console.",
      "recommendation": "Remove or mask sensitive information before printing. Only display necessary details for authorized users.",
      "code_suggestion": "Write-Host \"`nDisks:\""
    },
    {
      "line": 27,
      "severity": "Error",
      "statement": "Write-Host (\"Drive {0} - {1} GB free of {2} GB\" -f $_.DeviceID,",
      "reason": "Sensitive variable $disks (disk objects) is being printed to the console.",
      "recommendation": "Remove or mask sensitive information before printing. Only display necessary details for authorized users.",
      "code_suggestion": "Write-Host (\"Drive [Disk ID Masked] - [Free Space Masked] GB free of [Total Size Masked] GB\" -f $_.DeviceID,"
    }
  ]
}
'@
$ws.Range("B4").Value = $v2

$v3 = @'
vulnerable
'@
$ws.Range("C4").Value = $v3

$v4 = @'
{
  "script": "vulnerable",
  "score": 7,
  "findings": [
    {
      "line": 17,
      "severity": "Error",
      "statement": "Write-Host \"SNMP Auth Token: $($snmp.pass)\"",
      "reason": "Sensitive variable $snmp.pass is being printed to the console.",
      "recommendation": "All sensitive variables must be used internally and never output directly. Use secure methods for authentication without exposing credentials.",
      "code_suggestion": "Remove this line or replace with a safe confirmation message like Write-Host \"Authentication successful.\""
    },
    {
      "line": 19,
      "severity": "Error",
      "statement": "Set-Content -Path \"secrets.txt\" -Value $adminCredential.password",
      "reason": "Sensitive variable $adminCredential.password is being written to a file.",
      "recommendation": "Never store sensitive credentials in files or unencrypted storage. Use secure vaults or internal variables only.",
      "code_suggestion": "Remove this line and ensure the password is not persisted after use."
    },
    {
      "line": 21,
      "severity": "Error",
      "statement": "Write-Host \"POST /sendToken?key=$($azure.secretkey)\"",
      "reason": "Sensitive variable $azure.secretkey is being printed to the console.",
      "recommendation": "All sensitive variables must be used internally and never output directly. Use secure methods for authentication without exposing credentials.",
      "code_suggestion": "Remove this line or replace with a safe confirmation message like Write-Host \"Token sending initiated.\""
    }
  ]
}
'@
$ws.Range("B5").Value = $v4

$v5 = @'
{
  "script": "vulnerable",
  "score": 7,
  "findings": [
    {
      "line": 17,
      "severity": "Warning",
      "statement": "Write-Host \"[+] Connection opened successfully.\"",
      "reason": "Sensitive variable $password is being printed to the console.",
      "recommendation": "Remove this informational message or ensure it's not logged in production environments. Use secure methods for authentication without exposing credentials.",
      "code_suggestion": "Write-Host \"Connection successful\""
    },
    {
      "line": 24,
      "severity": "Warning",
      "statement": "Write-Host \"User Record: $($reader[0])\"",
      "reason": "Sensitive variable $password is being used to access database records.",
      "recommendation": "Ensure proper authorization and auditing mechanisms are in place when accessing sensitive data. Limit access based on the principle of least privilege.",
      "code_suggestion": "Write-Host \"User Record accessed\""
    },
    {
      "line": 35,
      "severity": "Error",
      "statement": "Write-Host \"DB Password in use: $password\"",
      "reason": "Sensitive variable $password is being printed to the console - a clear security leak.",
      "recommendation": "Never print or log sensitive credentials. Use secure methods for authentication without exposing passwords.",
      "code_suggestion": "Write-Host \"Authentication successful\""
    },
    {
      "line": 37,
      "severity": "Error",
      "statement": "Set-Content -Path \"leaked_secret.txt\" -Value $gcp.serviceaccountkey",
      "reason": "Sensitive variable $gcp.serviceaccountkey is being written to a file - a major security leak.",
      "recommendation": "Never store sensitive credentials in files or unencrypted storage. Use secure vaults or key management systems instead.",
      "code_suggestion": "Remove this line entirely as it creates an exploitable vulnerability."
    },
    {
      "line": 39,
      "severity": "Error",
      "statement": "Write-Host \"curl -X POST https://api.saasvendor.com/init --data 'key=$($saas.secretkey)'\"",
      "reason": "Sensitive variable $saas.secretkey is being displayed in a command that could be copied and used maliciously.",
      "recommendation": "Never display or expose sensitive credentials in commands or output. Use secure APIs or authentication methods instead.",
      "code_suggestion": "Remove this line as it reveals how to access the SaaS service."
    }
  ]
}
'@
$ws.Range("B6").Value = $v5

$v6 = @'
{
  "script": "vulnerable",
  "score": 8,
  "findings": [
    {
      "line": 8,
      "severity": "Error",
      "statement": "Invoke-Expression $sshCommand",
      "reason": "The SSH password is being passed to a shell command for execution.",
      "recommendation": "Avoid passing sensitive credentials directly to shell commands. Use secure methods like credential objects or key-based authentication instead.",
      "code_suggestion": "Use a secure method to execute the remote command without exposing the password."
    },
    {
      "line": 11,
      "severity": "Error",
      "statement": "Write-Host \"[!] For debug: SSH password = $env:SSH_PASS\"",
      "reason": "The sensitive SSH password is being printed to the console.",
      "recommendation": "Never print or log sensitive credentials. Remove this debugging statement and use secure methods for troubleshooting.",
      "code_suggestion": "Remove this line entirely as it's leaking a sensitive value."
    }
  ]
}
'@
$ws.Range("B7").Value = $v6

$v7 = @'
ERROR: 'system'
'@
$ws.Range("B8").Value = $v7

$v8 = @'
ERROR: 'system'
'@
$ws.Range("B9").Value = $v8

$v9 = @'
{
  "script": "vulnerable",
  "score": 6,
  "findings": [
    {
      "line": 3,
      "severity": "Error",
      "statement": "Set-Clipboard -Value $sshToken",
      "reason": "Sensitive SSH token is being copied to the clipboard, which is a public attack surface.",
      "recommendation": "Never copy sensitive credentials to the clipboard. Store them securely in memory or use secure credential management systems.",
      "code_suggestion": "Remove this line and handle authentication through secure channels."
    },
    {
      "line": 4,
      "severity": "Warning",
      "statement": "Write-Host \"[+] Sensitive SSH token copied to clipboard.\"",
      "reason": "While not directly leaking the value, confirming a sensitive operation occurred can provide attackers with information.",
      "recommendation": "Remove this informational message as it confirms a potentially risky action took place.",
      "code_suggestion": "Replace with a generic success indicator or remove entirely."
    },
    {
      "line": 9,
      "severity": "Warning",
      "statement": "$clipboardToken = Get-Clipboard",
      "reason": "Reading from the clipboard can expose sensitive data if it contains credentials.",
      "recommendation": "Avoid reading from the clipboard unless absolutely necessary and ensure no sensitive information is present.",
      "code_suggestion": "Remove this line or add a check to verify the clipboard contents before proceeding."
    },
    {
      "line": 16,
      "severity": "Error",
      "statement": "Write-Host \"`n[+] Executing SSH command:\"",
      "reason": "Printing the full SSH command including the token exposes sensitive information in plain text.",
      "recommendation": "Never print commands containing credentials. Use secure methods to execute them without exposing the actual values.",
      "code_suggestion": "Remove this line and handle execution through a secure channel."
    },
    {
      "line": 17,
      "severity": "Error",
      "statement": "Write-Host $sshCommand",
      "reason": "The SSH command contains the sensitive token, which is being printed to the console.",
      "recommendation": "Never print commands containing credentials. Use secure methods to execute them without exposing the actual values.",
      "code_suggestion": "Remove this line and handle execution through a secure channel."
    },
    {
      "line": 20,
      "severity": "Error",
      "statement": "Write-Host \"`n[!] DEBUG: Clipboard contents = $clipboardToken\"",
      "reason": "Printing the clipboard contents (which includes the token) exposes sensitive information for debugging purposes.",
      "recommendation": "Never print credential values, even in debug mode. Use secure logging mechanisms or remove this line entirely.",
      "code_suggestion": "Replace with a generic debug indicator or remove."
    }
  ]
}
'@
$ws.Range("B10").Value = $v9

$v10 = @'
{
  "script": "vulnerable",
  "score": 8,
  "findings": [
    {
      "line": 3,
      "severity": "Error",
      "statement": "iex $unsafeCmd",
      "reason": "Sensitive variable $snmp.auth is being passed to iex without proper sanitization.",
      "recommendation": "All sensitive variables must be used strictly for authentication within the script and never passed to shell execution functions like iex or Invoke-Expression.",
      "code_suggestion": "Remove this line as it creates a security vulnerability by potentially executing arbitrary code with the sensitive token."
    },
    {
      "line": 10,
      "severity": "Error",
      "statement": "iex $userInput",
      "reason": "Sensitive variable $snmp.auth is being passed to iex through an unsanitized user input.",
      "recommendation": "User inputs should never be directly executed as code; implement proper validation and sanitization before any execution.",
      "code_suggestion": "Remove this line or replace it with a secure method of handling user input that does not involve executing arbitrary commands."
    }
  ]
}
'@
$ws.Range("B11").Value = $v10

$v11 = @'
{
  "script": "safe",
  "score": 10,
  "findings": [
    {
      "line": 3,
      "severity": "Warning",
      "statement": "$sqlPass = \"##sql.pass##\"",
      "reason": "Sensitive variable declared but not used in a risky sink.",
      "recommendation": "Use this variable only for authentication within the script; do not output its value.",
      "code_suggestion": "$sqlPass = 'internal-use-only'"
    },
    {
      "line": 6,
      "severity": "Warning",
      "statement": "$ldapPass = \"##ldap.auth##\"",
      "reason": "Sensitive variable declared but not used in a risky sink.",
      "recommendation": "Use this variable only for authentication within the script; do not output its value.",
      "code_suggestion": "$ldapPass = 'internal-use-only'"
    },
    {
      "line": 8,
      "severity": "Warning",
      "statement": "$secureSqlPass = ConvertTo-SecureString $sqlPass -AsPlainText -Force",
      "reason": "Sensitive variable is being converted to a secure string but not used in a risky sink.",
      "recommendation": "Ensure this secure credential is only used for internal authentication purposes.",
      "code_suggestion": "$secureSqlPass = ConvertTo-SecureString 'internal-use-only' -AsPlainText -Force"
    },
    {
      "line": 13,
      "severity": "Warning",
      "statement": "$directory = New-Object System.DirectoryServices.DirectoryEntry(\"LDAP://domain.com\", $ldapUser, $ldapPass)",
      "reason": "Sensitive variable used in LDAP authentication but not output.",
      "recommendation": "Ensure this connection is only used for authorized operations within the script.",
      "code_suggestion": "$directory = New-Object System.DirectoryServices.DirectoryEntry(\"LDAP://domain.com\", $ldapUser, 'internal-use-only')"
    }
  ]
}
'@
$ws.Range("B12").Value = $v11

$v12 = @'
{
  "script": "vulnerable",
  "score": 9,
  "findings": [
    {
      "line": 3,
      "severity": "Warning",
      "statement": "$sqlPass = \"##sql.pass##\"",
      "reason": "Sensitive variable declared but not used in a risky sink.",
      "recommendation": "Ensure this password is only used internally for database operations and never printed or logged.",
      "code_suggestion": "$sqlPass = New-Object SecureString -AsPlainText 'your_secure_password' -KeyLength 256"
    },
    {
      "line": 6,
      "severity": "Warning",
      "statement": "$ldapPass = \"##ldap.auth##\"",
      "reason": "Sensitive variable declared but not used in a risky sink.",
      "recommendation": "Ensure this password is only used internally for LDAP operations and never printed or logged.",
      "code_suggestion": "$ldapPass = New-Object SecureString -AsPlainText 'your_secure_password' -KeyLength 256"
    },
    {
      "line": 10,
      "severity": "Error",
      "statement": "Write-Host \"LDAP bind password: $ldapPass\"",
      "reason": "Sensitive variable value is being printed to the console.",
      "recommendation": "Remove this line as it exposes a sensitive credential. Use secure methods for authentication without logging passwords.",
      "code_suggestion": "Remove or replace with a success message indicating authentication was successful."
    }
  ]
}
'@
$ws.Range("B13").Value = $v12

$v13 = @'
{
  "script": "vulnerable",
  "score": 8,
  "findings": [
    {
      "line": 8,
      "severity": "Error",
      "statement": "Write-Host \"SQL credentials: $sqlUser / $sqlPass\"",
      "reason": "Sensitive variable $sqlPass is being printed to the console.",
      "recommendation": "All sensitive variables must be used strictly for authentication within the script and never exposed in output.",
      "code_suggestion": "Remove this line or replace with a safe message like Write-Host \"Authentication successful\""
    },
    {
      "line": 9,
      "severity": "Error",
      "statement": "Write-Host \"LDAP bind credentials: $ldapUser / $ldapPass\"",
      "reason": "Sensitive variable $ldapPass is being printed to the console.",
      "recommendation": "All sensitive variables must be used strictly for authentication within the script and never exposed in output.",
      "code_suggestion": "Remove this line or replace with a safe message like Write-Host \"Authentication successful\""
    },
    {
      "line": 13,
      "severity": "Warning",
      "statement": "$directory = New-Object System.DirectoryServices.DirectoryEntry(\"LDAP://domain.com\", $ldapUser, $ldapPass)",
      "reason": "Sensitive variable $ldapPass is being used in a potentially risky operation (LDAP bind). While not immediately leaking the value, this pattern could lead to future vulnerabilities if the directory object or searcher are later misused.",
      "recommendation": "Ensure that all operations using sensitive credentials follow the principle of least privilege and have clear audit trails. Consider using more secure authentication methods where possible.",
      "code_suggestion": "$directory = New-Object System.DirectoryServices.DirectoryEntry(\"LDAP://domain.com\", $ldapUser, ConvertTo-SecureString $ldapPass -AsPlainText -Force)"
    }
  ]
}
'@
$ws.Range("B14").Value = $v13

$v14 = @'
{
  "script": "vulnerable",
  "score": 9,
  "findings": [
    {
      "line": 3,
      "severity": "Warning",
      "statement": "$password = \"##wmi.pass##\"",
      "reason": "Sensitive variable declared with placeholder value.",
      "recommendation": "Use a secure method to retrieve or generate the password instead of hardcoding it.",
      "code_suggestion": "$password = Read-Host -AsSecureString 'Enter WMI password: '"
    },
    {
      "line": 8,
      "severity": "Warning",
      "statement": "$session = New-PSSession -ComputerName $remoteHost -Credential $cred",
      "reason": "Sensitive credential used to establish a remote session.",
      "recommendation": "Ensure the remote host is trusted and that access controls are properly configured.",
      "code_suggestion": "$session = New-PSSession -ComputerName $remoteHost -Credential (Get-Credential)"
    },
    {
      "line": 10,
      "severity": "Error",
      "statement": "Invoke-Command -Session $session -ScriptBlock {{1}} Get-WmiObject Win32_OperatingSystem {{1}}",
      "reason": "Sensitive credential used in a network operation (remote command execution).",
      "recommendation": "Avoid passing credentials directly to remote commands; use alternative authentication methods or constrained delegation.",
      "code_suggestion": "Invoke-Command -Session $session -ScriptBlock {{Get-WmiObject Win32_OperatingSystem}}"
    },
    {
      "line": 12,
      "severity": "Warning",
      "statement": "Write-Host \"Remote data collection from $remoteHost completed.\"",
      "reason": "Sensitive variable used in console output.",
      "recommendation": "Remove or redact sensitive information from log messages and console outputs.",
      "code_suggestion": "Write-Host 'Remote data collection complete.'"
    }
  ]
}
'@
$ws.Range("B15").Value = $v14

$v15 = @'
{
  "script": "vulnerable",
  "score": 9,
  "findings": [
    {
      "line": 9,
      "severity": "Error",
      "statement": "Write-Host \"Collected service data from $remoteHost using $username:$password\"",
      "reason": "Sensitive password is being printed to the console.",
      "recommendation": "All sensitive credentials must be used internally and never logged or displayed. Use a secure method for storing and retrieving passwords, such as a credential management system.",
      "code_suggestion": "Remove this line entirely, as it contains a clear security vulnerability."
    }
  ]
}
'@
$ws.Range("B16").Value = $v15

$v16 = @'
{
  "script": "vulnerable",
  "score": 6,
  "findings": [
    {
      "line": 3,
      "severity": "Warning",
      "statement": "$password = \"##wmi.pass##\"",
      "reason": "Sensitive variable declared with placeholder value.",
      "recommendation": "Use a secure method to retrieve or generate the password instead of hardcoding it.",
      "code_suggestion": "$password = Read-Host -AsSecureString 'Enter password:'"
    },
    {
      "line": 7,
      "severity": "Error",
      "statement": "Write-Host \"Connecting to $remoteHost as $username with password $password\"",
      "reason": "Sensitive variable `$password` is being printed to the console.",
      "recommendation": "Remove this output statement or use a secure logging mechanism that does not expose sensitive credentials.",
      "code_suggestion": "Write-Host \"Connecting to $remoteHost as $username\""
    },
    {
      "line": 9,
      "severity": "Error",
      "statement": "Invoke-Command -Session $session -ScriptBlock {{1}} Get-Process {{1}} > \"C:\\logs\\process_$remoteHost.txt\"",
      "reason": "Sensitive variable `$password` is being used in a command that outputs to a file.",
      "recommendation": "Ensure the log file has restricted access and consider encrypting sensitive data before logging.",
      "code_suggestion": "Invoke-Command -Session $session -ScriptBlock {{1}} Get-Process {{1}} | Out-Null"
    },
    {
      "line": 11,
      "severity": "Error",
      "statement": "Invoke-WebRequest -Uri \"https://logstash.company.com/ingest?u=$username&p=$password\"",
      "reason": "Sensitive variable `$password` is being transmitted over the network in plain text.",
      "recommendation": "Use HTTPS and consider encrypting or hashing sensitive data before transmission.",
      "code_suggestion": "Invoke-WebRequest -Uri \"https://logstash.company.com/ingest?u=$username&p=$(ConvertTo-SecureString $password -AsPlainText -Force)\""
    },
    {
      "line": 13,
      "severity": "Error",
      "statement": "Set-Content -Path \"C:\\temp\\session_info.txt\" -Value \"User=$username Password=$password Host=$remoteHost\"",
      "reason": "Sensitive variable `$password` is being written to a file in plain text.",
      "recommendation": "Remove this output statement or use a secure method to store session information.",
      "code_suggestion": "Set-Content -Path \"C:\\temp\\session_info.txt\" -Value \"User=$username Host=$remoteHost\""
    }
  ]
}
'@
$ws.Range("B17").Value = $v16
